$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the existing tag table (A1:E5) before rewriting it with the
# new layout produced by the updated query-generation algorithm.
$ws.Range("A1:F5").ClearContents()

# New column layout: each column is a tag-group header (row 1) followed
# by its member tags. A new "intersection_operators" group (with the
# single tag "has posted") has been inserted as the second column,
# pushing the previous join/clause/collateral/legal_entity groups one
# column to the right.
$columns = @{
    "A" = @("initial_operators", "select", "filter")
    "B" = @("intersection_operators", "has posted")
    "C" = @("join_operators", "and", "but not")
    "D" = @("clause_operators", "where", "in")
    "E" = @("collateral", "inrcash", "uscorp", "gbpcash", "usmuni")
    "F" = @("legal_entity", "cgmi", "cgml", "cgma")
}

foreach ($col in $columns.Keys) {
    $values = $columns[$col]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 1
        $ws.Range("$col$row").Value = $values[$i]
    }
}

$ws.Range("D3").Select()
